$wb = $excel.ActiveWorkbook

# --- ALC row 40 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2063.3333
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2063.3333
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2063.3333
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2413.3333

# --- ALC row 118 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 959.1539
$ws.Range("I118").Value = 574.3333
$ws.Range("J118").Value = 1825
$ws.Range("K118").Value = 1722.9999
$ws.Range("L118").Value = 5475
$ws.Range("M118").Value = -65.99990000000003
$ws.Range("N118").Value = -8789

# --- ALC row 132 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2441502
$ws.Range("I132").Value = 2415.8667
$ws.Range("J132").Value = 9093555
$ws.Range("K132").Value = 7247.6001
$ws.Range("L132").Value = 27280665
$ws.Range("M132").Value = -4717.6001
$ws.Range("N132").Value = -27285725

# --- ALC row 133 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 46548.234
$ws.Range("I133").Value = 25800
$ws.Range("J133").Value = 49314.668
$ws.Range("K133").Value = 25800
$ws.Range("L133").Value = 49314.668
$ws.Range("M133").Value = -20740
$ws.Range("N133").Value = -59434.668

# --- ALC row 138 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4039.8386
$ws.Range("I138").Value = 1707.2354
$ws.Range("J138").Value = 4561.6055
$ws.Range("K138").Value = 5121.706200000001
$ws.Range("L138").Value = 13684.8165
$ws.Range("M138").Value = 18.29379999999946
$ws.Range("N138").Value = -23964.8165

# --- ARM row 32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18431.623
$ws.Range("I32").Value = 13394.345
$ws.Range("J32").Value = 56840.875
$ws.Range("K32").Value = 13394.345
$ws.Range("L32").Value = 56840.875
$ws.Range("M32").Value = -13107.345
$ws.Range("N32").Value = -57414.875

# --- ARM row 37 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 6914.2
$ws.Range("I37").Value = 2534
$ws.Range("J37").Value = 8009.25
$ws.Range("K37").Value = 2534
$ws.Range("L37").Value = 8009.25
$ws.Range("M37").Value = -2261
$ws.Range("N37").Value = -8555.25

# --- ARM row 44 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 21436.074
$ws.Range("I44").Value = 5000
$ws.Range("J44").Value = 22068.23
$ws.Range("K44").Value = 5000
$ws.Range("L44").Value = 22068.23
$ws.Range("M44").Value = -4512
$ws.Range("N44").Value = -23044.23

# --- ARM row 55 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 25064.297
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 25064.297
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 25064.297
$ws.Range("N55").Value = -25694.297

# --- ARM row 80 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 24037.691
$ws.Range("I80").Value = 500
$ws.Range("J80").Value = 25999.166
$ws.Range("K80").Value = 500
$ws.Range("L80").Value = 25999.166
$ws.Range("M80").Value = 498
$ws.Range("N80").Value = -27995.166

# --- ARM row 83 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 24037.691
$ws.Range("I83").Value = 500
$ws.Range("J83").Value = 25999.166
$ws.Range("K83").Value = 1500
$ws.Range("L83").Value = 77997.49800000001
$ws.Range("M83").Value = 3492
$ws.Range("N83").Value = -87981.49800000001

# --- ARM row 102 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 10000
$ws.Range("I102").Value = 10000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 10000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -8378
$ws.Range("N102").ClearContents()

# --- ARM row 122 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1418.1818
$ws.Range("I122").Value = 1400
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4200
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1750
$ws.Range("N122").Value = -9400

# --- ARM row 140 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 70952.664
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 70952.664
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 70952.664
$ws.Range("N140").Value = -81312.664

# --- BSM row 141 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 54906.832
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 54906.832
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 54906.832
$ws.Range("N141").Value = -65266.832

# --- CRP row 31 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3548.5715
$ws.Range("I31").Value = 3227.7368
$ws.Range("J31").Value = 4225.8887
$ws.Range("K31").Value = 3227.7368
$ws.Range("L31").Value = 4225.8887
$ws.Range("M31").Value = -2932.7368
$ws.Range("N31").Value = -4815.8887

# --- CRP row 34 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3548.5715
$ws.Range("I34").Value = 3227.7368
$ws.Range("J34").Value = 4225.8887
$ws.Range("K34").Value = 3227.7368
$ws.Range("L34").Value = 4225.8887
$ws.Range("M34").Value = -3025.7368
$ws.Range("N34").Value = -4629.8887

# --- CRP row 112 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 24980
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 24980
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 24980
$ws.Range("N112").Value = -27934

# --- CRP row 132 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1584.6
$ws.Range("I132").Value = 1152.5294
$ws.Range("J132").Value = 2502.75
$ws.Range("K132").Value = 3457.5882
$ws.Range("L132").Value = 7508.25
$ws.Range("M132").Value = -927.5881999999997
$ws.Range("N132").Value = -12568.25

# --- CUL row 117 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 2418.2
$ws.Range("I117").Value = 340.16666
$ws.Range("J117").Value = 3803.5557
$ws.Range("K117").Value = 1020.49998
$ws.Range("L117").Value = 11410.6671
$ws.Range("M117").Value = 2421.50002
$ws.Range("N117").Value = -18294.6671

# --- CUL row 121 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 875.4286
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 875.4286
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 2626.2858
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -5246.2858

# --- CUL row 122 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 22510.404
$ws.Range("I122").Value = 596.2
$ws.Range("J122").Value = 25119.238
$ws.Range("K122").Value = 5365.8
$ws.Range("L122").Value = 226073.142
$ws.Range("M122").Value = -2915.8
$ws.Range("N122").Value = -230973.142

# --- CUL row 131 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 20081326
$ws.Range("I131").Value = 62625480
$ws.Range("J131").Value = 60548.293
$ws.Range("K131").Value = 187876440
$ws.Range("L131").Value = 181644.879
$ws.Range("M131").Value = -187871400
$ws.Range("N131").Value = -191724.879

# --- GSM row 70 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4425.6665
$ws.Range("I70").Value = 4014.1428
$ws.Range("J70").Value = 5001.8
$ws.Range("K70").Value = 4014.1428
$ws.Range("L70").Value = 5001.8
$ws.Range("M70").Value = -3744.1428
$ws.Range("N70").Value = -5541.8

# --- GSM row 73 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4425.6665
$ws.Range("I73").Value = 4014.1428
$ws.Range("J73").Value = 5001.8
$ws.Range("K73").Value = 4014.1428
$ws.Range("L73").Value = 5001.8
$ws.Range("M73").Value = -3078.1428
$ws.Range("N73").Value = -6873.8

# --- GSM row 122 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1097404.6
$ws.Range("I122").Value = 1880597.1
$ws.Range("J122").Value = 935.2
$ws.Range("K122").Value = 5641791.300000001
$ws.Range("L122").Value = 2805.6
$ws.Range("M122").Value = -5639341.300000001
$ws.Range("N122").Value = -7705.6

# --- GSM row 123 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 50093.145
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 50093.145
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 50093.145
$ws.Range("N123").Value = -54993.145

# --- LTW row 22 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 735
$ws.Range("I22").Value = 533.3333
$ws.Range("J22").Value = 790
$ws.Range("K22").Value = 533.3333
$ws.Range("L22").Value = 790
$ws.Range("M22").Value = -238.3333
$ws.Range("N22").Value = -1380

# --- LTW row 27 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 735
$ws.Range("I27").Value = 533.3333
$ws.Range("J27").Value = 790
$ws.Range("K27").Value = 533.3333
$ws.Range("L27").Value = 790
$ws.Range("M27").Value = -426.3333
$ws.Range("N27").Value = -1004

# --- LTW row 46 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 50001364
$ws.Range("I46").Value = 50001364
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 50001364
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -50001176
$ws.Range("N46").ClearContents()

# --- LTW row 55 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 187.92308
$ws.Range("I55").Value = 77.57143000000001
$ws.Range("J55").Value = 316.66666
$ws.Range("K55").Value = 77.57143000000001
$ws.Range("L55").Value = 316.66666
$ws.Range("M55").Value = 95.42856999999999
$ws.Range("N55").Value = -662.66666
